# Insert a new column before column A, shifting the existing
# Gene Name / Alias1 / Alias2 / Alias3 table one column to the right
# (A-D -> B-E), then fill the new column A (rows 2-103) with a
# 0-based row index, styled like the header row (bold, bordered,
# centered) by copying the header's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns A:D -> B:E
$ws.Columns("A").Insert()

# Fill the new column A with a 0-based index for every data row
$lastRow = 103
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Give the new index cells the same look as the header row
# (bold font, thin border, centered/top aligned) by copying the
# header cell's formatting onto them.
$ws.Range("B1").Copy()
$ws.Range("A2:A103").PasteSpecial(-4122)
$excel.CutCopyMode = $false
